$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.3980987966060638
$ws.Range("B1").Value = 1.995837450027466
$ws.Range("C1").Value = 4.886344909667969
$ws.Range("D1").Value = 1.736964821815491
$ws.Range("E1").Value = 0.8746339678764343
